$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Wrap the SPROuT logo picture (currently a lone "First Paragraph"
#    paragraph) in a one-cell table, and add a second (currently empty)
#    "Image Caption" placeholder paragraph below the picture.
# ---------------------------------------------------------------------
$logoParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "First Paragraph" -and $p.Range.InlineShapes.Count -gt 0) {
        $logoParagraph = $p
        break
    }
}

$tableXml = '<w:tbl xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"><w:tblPr><w:tblStyle w:val="Table" /><w:tblW w:type="pct" w:w="5000" /><w:tblLook w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0" w:val="0000" /><w:jc w:val="start" /></w:tblPr><w:tblGrid><w:gridCol w:w="7920" /></w:tblGrid><w:tr><w:tc><w:tcPr /><w:p><w:pPr><w:jc w:val="center" /></w:pPr><w:r><w:drawing><wp:inline><wp:extent cx="2857500" cy="2857500" /><wp:effectExtent b="0" l="0" r="0" t="0" /><wp:docPr descr="" title="" id="21" name="Picture" /><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="sprout_logo_blue.jpg" id="22" name="Picture" /><pic:cNvPicPr><a:picLocks noChangeArrowheads="1" noChangeAspect="1" /></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId20" /><a:stretch><a:fillRect /></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0" /><a:ext cx="2857500" cy="2857500" /></a:xfrm><a:prstGeom prst="rect"><a:avLst /></a:prstGeom><a:noFill /><a:ln w="9525"><a:noFill /><a:headEnd /><a:tailEnd /></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p><w:p><w:pPr><w:jc w:val="center" /></w:pPr><w:pPr><w:jc w:val="start" /><w:spacing w:before="200" /><w:pStyle w:val="ImageCaption" /></w:pPr></w:p></w:tc></w:tr></w:tbl>'

$logoParagraph.Range.InsertXML($tableXml)

# ---------------------------------------------------------------------
# 2) Styles: drop the unused "Abstract Title" style, fold its intent
#    into "Abstract" (keep "Abstract", but bump its space-before from
#    100 twips/5pt to 300 twips/15pt so it matches what the title used
#    to add), and drop the unused "Footnote Block Text" style.
# ---------------------------------------------------------------------
$d.Styles.Item("AbstractTitle").Delete()
$d.Styles.Item("Abstract").ParagraphFormat.SpaceBefore = 15
$d.Styles.Item("FootnoteBlockText").Delete()
